# Updated cryptos list on Thu Dec 21 16:27:10 UTC 2023 with GitHub Actions
#
# Refreshes the coin Price (column D) and Volume(1h) (column E) figures for
# the rows that moved, and fixes the Uniswap/Litecoin row ordering (rows 20
# and 21 had their Coin/Link/Price/Volume data swapped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> column/value pairs to write (only the cells that actually changed)
$updates = [ordered]@{
    2  = @{ D = '43.778.76'; E = '  -0.55%  ' }
    3  = @{ D = '2.227.28';  E = '  -0.71%  ' }
    4  = @{ E = '  +0.00%  ' }
    5  = @{ D = '271.69';    E = '  +5.12%  ' }
    6  = @{ D = '89.26';     E = '  +11.58%  ' }
    7  = @{ D = '0.621';     E = '  -0.71%  ' }
    8  = @{ E = '  +0.00%  ' }
    9  = @{ D = '0.605';     E = '  +0.50%  ' }
    10 = @{ D = '45.82';     E = '  +5.77%  ' }
    11 = @{ D = '0.0918';    E = '  -1.20%  ' }
    12 = @{ D = '7.74';      E = '  +9.30%  ' }
    13 = @{ D = '0.105';     E = '  +1.34%  ' }
    14 = @{ D = '2.563.82';  E = '  -0.43%  ' }
    15 = @{ D = '15.01';     E = '  +2.27%  ' }
    16 = @{ D = '2.213.23';  E = '  -2.66%  ' }
    17 = @{ D = '0.793';     E = '  +0.29%  ' }
    18 = @{ D = '43.739.44'; E = '  -0.41%  ' }
    19 = @{ D = '0.0000103'; E = '  -1.27%  ' }
    20 = @{ B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '70.25'; E = '  -1.63%  ' }
    21 = @{ B = 'Uniswap';  C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni';    D = '5.97';  E = '  -1.37%  ' }
    22 = @{ D = '2.36' }
    23 = @{ D = '232.18';   E = '  -0.50%  ' }
    24 = @{ D = '8.60';     E = '  -8.06%  ' }
    26 = @{ D = '2.52';     E = '  +13.23%  ' }
    27 = @{ D = '10.95';    E = '  +0.98%  ' }
    28 = @{ D = '3.56';     E = '  +5.60%  ' }
    29 = @{ D = '2.28';     E = '  +3.03%  ' }
    30 = @{ D = '38.80';    E = '  -5.03%  ' }
    31 = @{ D = '173.07';   E = '  +0.17%  ' }
    32 = @{ D = '0.0912';   E = '  +2.54%  ' }
    33 = @{ D = '20.70';    E = '  +0.45%  ' }
    34 = @{ E = '  +0.74%  ' }
    35 = @{ E = '  -0.08%  ' }
    36 = @{ E = '  -1.01%  ' }
    37 = @{ D = '0.0353';   E = '  -3.18%  ' }
    38 = @{ D = '4.25';     E = '  -6.27%  ' }
    39 = @{ D = '3.46';     E = '  +17.17%  ' }
    40 = @{ D = '2.16';     E = '  +0.83%  ' }
    41 = @{ D = '12.34';    E = '  -4.68%  ' }
    42 = @{ E = '  +4.99%  ' }
    43 = @{ D = '63.37';    E = '  +0.64%  ' }
    44 = @{ D = '5.38';     E = '  -2.86%  ' }
    45 = @{ D = '8.50';     E = '  -0.23%  ' }
    46 = @{ D = '0.0985';   E = '  -0.10%  ' }
    47 = @{ D = '99.95';    E = '  -4.04%  ' }
    48 = @{ D = '1.14';     E = '  +1.48%  ' }
    49 = @{ D = '1.18';     E = '  +2.68%  ' }
    50 = @{ D = '0.433';    E = '  -3.39%  ' }
    51 = @{ D = '1.47';     E = '  -3.64%  ' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in @('B', 'C', 'D', 'E')) {
        if ($cols.Contains($col)) {
            $address = "$col$row"
            $value = $cols[$col]
            $cell = $ws.Range($address)
            if ($col -eq 'D') {
                # Column D holds prices like "43.778.76" / "0.0918" which are
                # stored as plain text in the sheet. Force text formatting
                # before the assignment so Excel doesn't silently reinterpret
                # numeric-looking strings (e.g. "45.82", "0.621") as numbers,
                # then restore the default (unstyled) look of the data rows.
                $cell.NumberFormat = "@"
                $cell.Value = $value
                $cell.Style = "Normal"
            } else {
                $cell.Value = $value
            }
        }
    }
}
